# Update Excel files for new format
# Adds a new "Setup" sheet (listing the header values used by the
# ConflictMatrix together with a constant 1.5 column) in front of the
# existing matrix sheet, which is renamed from "Sheet1" to "ConflictMatrix".

$wb = $excel.ActiveWorkbook

# --- rename the existing sheet first --------------------------------------
$matrix = $wb.Worksheets.Item("Sheet1")
$matrix.Name = "ConflictMatrix"

# --- insert the new "Setup" sheet in front of it ---------------------------
$setup = $wb.Worksheets.Add()
$setup.Name = "Setup"

# Column A: the conflicting-signal-group id used as the ConflictMatrix's
# row/column headers. Column B: a constant 1.5 (default clearance time).
$ids = @(
    1.1, 2.1, 5.1, 6.1, 7.1, 8.1, 9.1, 10.1, 11.1, 12.1,
    86.1, 35.1, 26.1, 36.2, 88.1, 37.2, 28.1, 38.2, 31.2, 22,
    32.2, 35.2, 36.1, 37.1, 38.1, 31.1, 32.1
)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 1
    $setup.Cells.Item($row, 1).Value = $ids[$i]
    $setup.Cells.Item($row, 2).Value = 1.5
}
